$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.08897959183673469
$ws.Range("A3").Value = 0.290204081632653
$ws.Range("A4").Value = 0.6032653061224489
$ws.Range("A5").Value = 0.8285714285714285
$ws.Range("A6").Value = 0.9306122448979591
$ws.Range("A7").Value = 0.966938775510204
$ws.Range("A8").Value = 0.986938775510204
$ws.Range("A9").Value = 0.9946938775510203
$ws.Range("A10").Value = 0.996734693877551
$ws.Range("A11").Value = 0.996734693877551
$ws.Range("A13").Value = 0.9971428571428571
$ws.Range("A14").Value = 0.9975510204081632
$ws.Range("A16").Value = 0.9975510204081632
